$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Begin working on numbers less than 1: new working columns (H, I) for
# exponent / mantissa extraction from the floating point registers.
#
# Cells are written in shared-string-creation order (exponent, bits in
# register, result of shift, nth register mantissa check, more than 22
# bits read in?, final mantissa) so the sharedStrings.xml table comes out
# in the same order Excel itself produced it in.
$ws.Range("H7").Value = "exponent"
$ws.Range("H9").Value = "bits in register"
$ws.Range("H10").Value = "result of shift"
$ws.Range("I1").Value = "nth register mantissa check"
$ws.Range("H11").Value = "more than 22 bits read in?"
$ws.Range("I14").Value = "final mantissa"

# Give the two new columns an explicit (best-fit-like) width, matching the
# pattern already used for columns B-G.
$ws.Columns.Item(8).ColumnWidth = 10.92
$ws.Columns.Item(9).ColumnWidth = 24.59

# Reflect the work happening further to the right of the sheet: the user
# ends up with I16 selected.
$ws.Range("I16").Select()
